# Automatic update of files.
# Bump the "Förändrad" (changed) date in column C, rows 2-9, by one day
# (from serial 46060 to 46061), keeping existing cell formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 9; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 46060) {
        $cell.Value2 = 46061
    }
}
